# Updates the "cryptos" price list (Coin / Link / Price / Volume(1h)) to the
# latest scraped snapshot. Most rows only get fresh Price / Volume(1h)
# figures, but a few coins were re-ranked and swapped positions (rows
# 27/28 and 50/51 each trade a whole row's Coin/Link/Price/Volume), and
# row 50/51's old "BabyDogeCoin" slot is replaced by "Algorand" while a
# brand new "EnergySwap" listing takes the old last row.
#
# Source data values like "1.00", "212.42" or "7.80" are display strings,
# not numbers (trailing zeros / thousands-style dots must survive as
# typed). Excel's COM layer auto-detects such strings as numeric on
# assignment, so for any replacement value that parses as a plain number
# we force the cell to Text format first, then restore the cell's style
# to Normal afterwards so only the *value* changes (no stray number
# format is left behind on the cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "27.699.62"; ForceText = $false },
    @{ Cell = "E2"; Value = "  +0.93%  "; ForceText = $false },
    @{ Cell = "D3"; Value = "1.627.35"; ForceText = $false },
    @{ Cell = "E3"; Value = "  -0.56%  "; ForceText = $false },
    @{ Cell = "E4"; Value = "  +0.27%  "; ForceText = $false },
    @{ Cell = "D5"; Value = "212.42"; ForceText = $true },
    @{ Cell = "E5"; Value = "  +0.06%  "; ForceText = $false },
    @{ Cell = "E6"; Value = "  -1.78%  "; ForceText = $false },
    @{ Cell = "E7"; Value = "  +0.24%  "; ForceText = $false },
    @{ Cell = "D8"; Value = "22.95"; ForceText = $true },
    @{ Cell = "E8"; Value = "  +0.37%  "; ForceText = $false },
    @{ Cell = "E9"; Value = "  +2.14%  "; ForceText = $false },
    @{ Cell = "D10"; Value = "0.0612"; ForceText = $true },
    @{ Cell = "E10"; Value = "  +0.29%  "; ForceText = $false },
    @{ Cell = "D11"; Value = "0.0889"; ForceText = $true },
    @{ Cell = "E11"; Value = "  +0.30%  "; ForceText = $false },
    @{ Cell = "D12"; Value = "1.857.98"; ForceText = $false },
    @{ Cell = "E12"; Value = "  -0.58%  "; ForceText = $false },
    @{ Cell = "D13"; Value = "1.625.36"; ForceText = $false },
    @{ Cell = "E13"; Value = "  -0.65%  "; ForceText = $false },
    @{ Cell = "D14"; Value = "4.05"; ForceText = $true },
    @{ Cell = "E14"; Value = "  +0.82%  "; ForceText = $false },
    @{ Cell = "E15"; Value = "  -2.89%  "; ForceText = $false },
    @{ Cell = "E16"; Value = "  +0.52%  "; ForceText = $false },
    @{ Cell = "D17"; Value = "27.686.08"; ForceText = $false },
    @{ Cell = "E17"; Value = "  +0.76%  "; ForceText = $false },
    @{ Cell = "D18"; Value = "230.03"; ForceText = $true },
    @{ Cell = "E18"; Value = "  +0.79%  "; ForceText = $false },
    @{ Cell = "E19"; Value = "  -0.11%  "; ForceText = $false },
    @{ Cell = "D20"; Value = "7.59"; ForceText = $true },
    @{ Cell = "E20"; Value = "  -0.67%  "; ForceText = $false },
    @{ Cell = "E21"; Value = "  +0.37%  "; ForceText = $false },
    @{ Cell = "D22"; Value = "4.31"; ForceText = $true },
    @{ Cell = "E22"; Value = "  -0.04%  "; ForceText = $false },
    @{ Cell = "D23"; Value = "9.97"; ForceText = $true },
    @{ Cell = "E23"; Value = "  +2.11%  "; ForceText = $false },
    @{ Cell = "E24"; Value = "  +6.17%  "; ForceText = $false },
    @{ Cell = "D25"; Value = "149.76"; ForceText = $true },
    @{ Cell = "E25"; Value = "  +0.08%  "; ForceText = $false },
    @{ Cell = "E26"; Value = "  -1.03%  "; ForceText = $false },
    @{ Cell = "B27"; Value = "Stellar"; ForceText = $false },
    @{ Cell = "C27"; Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; ForceText = $false },
    @{ Cell = "D27"; Value = "0.111"; ForceText = $true },
    @{ Cell = "E27"; Value = "  -0.84%  "; ForceText = $false },
    @{ Cell = "B28"; Value = "BinanceUSD"; ForceText = $false },
    @{ Cell = "C28"; Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"; ForceText = $false },
    @{ Cell = "D28"; Value = "1.00"; ForceText = $true },
    @{ Cell = "E28"; Value = "  +0.24%  "; ForceText = $false },
    @{ Cell = "D29"; Value = "15.58"; ForceText = $true },
    @{ Cell = "E29"; Value = "  +0.31%  "; ForceText = $false },
    @{ Cell = "E30"; Value = "  -0.29%  "; ForceText = $false },
    @{ Cell = "D31"; Value = "0.0482"; ForceText = $true },
    @{ Cell = "E31"; Value = "  -1.09%  "; ForceText = $false },
    @{ Cell = "E32"; Value = "  +0.27%  "; ForceText = $false },
    @{ Cell = "D33"; Value = "1.462.64"; ForceText = $false },
    @{ Cell = "E33"; Value = "  +3.51%  "; ForceText = $false },
    @{ Cell = "D34"; Value = "3.09"; ForceText = $true },
    @{ Cell = "E34"; Value = "  -2.44%  "; ForceText = $false },
    @{ Cell = "E35"; Value = "  -2.43%  "; ForceText = $false },
    @{ Cell = "E36"; Value = "  +0.70%  "; ForceText = $false },
    @{ Cell = "E37"; Value = "  -0.45%  "; ForceText = $false },
    @{ Cell = "E38"; Value = "  -0.10%  "; ForceText = $false },
    @{ Cell = "D39"; Value = "0.871"; ForceText = $true },
    @{ Cell = "E39"; Value = "  -0.35%  "; ForceText = $false },
    @{ Cell = "E40"; Value = "  -2.16%  "; ForceText = $false },
    @{ Cell = "D41"; Value = "69.40"; ForceText = $true },
    @{ Cell = "E41"; Value = "  +7.16%  "; ForceText = $false },
    @{ Cell = "D42"; Value = "1.00"; ForceText = $true },
    @{ Cell = "E42"; Value = "  +0.34%  "; ForceText = $false },
    @{ Cell = "E43"; Value = "  -1.17%  "; ForceText = $false },
    @{ Cell = "E44"; Value = "  -1.31%  "; ForceText = $false },
    @{ Cell = "D45"; Value = "5.44"; ForceText = $true },
    @{ Cell = "E45"; Value = "  -1.26%  "; ForceText = $false },
    @{ Cell = "E46"; Value = "  -0.28%  "; ForceText = $false },
    @{ Cell = "D47"; Value = "1.768.38"; ForceText = $false },
    @{ Cell = "E47"; Value = "  -0.56%  "; ForceText = $false },
    @{ Cell = "E48"; Value = "  +2.43%  "; ForceText = $false },
    @{ Cell = "D49"; Value = "85.88"; ForceText = $true },
    @{ Cell = "E49"; Value = "  +0.12%  "; ForceText = $false },
    @{ Cell = "B50"; Value = "Algorand"; ForceText = $false },
    @{ Cell = "C50"; Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; ForceText = $false },
    @{ Cell = "D50"; Value = "0.0989"; ForceText = $true },
    @{ Cell = "E50"; Value = "  +0.23%  "; ForceText = $false },
    @{ Cell = "B51"; Value = "EnergySwap"; ForceText = $false },
    @{ Cell = "C51"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; ForceText = $false },
    @{ Cell = "D51"; Value = "7.80"; ForceText = $true },
    @{ Cell = "E51"; Value = "  +1.44%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
